$d = $word.ActiveDocument

# 1. Report the fish x beef x breed ANOVA's breed p-value to two decimal
#    places instead of the full float.
$ok1 = $d.Content.Find.Execute("0.0343720051497151", $true, $false, $false, $false, $false,
                                $true, 1, $false, "0.03", 2)
Write-Output "edit1 (p-value): $ok1"

# 2. Add the missing space between "significant" and "(all" in the
#    interactions sentence.
$ok2 = $d.Content.Find.Execute("significant(all", $true, $false, $false, $false, $false,
                                $true, 1, $false, "significant (all", 2)
Write-Output "edit2 (space before parenthesis): $ok2"

# 3. "difference" -> "differences" (there are two comparisons being
#    described, so the plural is correct).
$ok3 = $d.Content.Find.Execute("The difference between Manx", $true, $false, $false, $false, $false,
                                $true, 1, $false, "The differences between Manx", 2)
Write-Output "edit3 (difference -> differences): $ok3"
